{"js": "// This table-driven script swaps each three-digit-division answer\n// cell for its replacement, addressing cells positionally (row, col)\n// so it is immune to any text collisions between old/new values\n// (e.g. one cell's new answer happens to equal another cell's old\n// answer elsewhere in the table).\nconst edits = [\n  {\n    \"row\": 0,\n    \"col\": 0,\n    \"find\": \"793\u00f75=158, 3\",\n    \"replace\": \"986\u00f74=246, 2\"\n  },\n  {\n    \"row\": 0,\n    \"col\": 1,\n    \"find\": \"489\u00f74=122, 1\",\n    \"replace\": \"231\u00f77=33, 0\"\n  },\n  {\n    \"row\": 0,\n    \"col\": 2,\n    \"find\": \"371\u00f78=46, 3\",\n    \"replace\": \"644\u00f77=92, 0\"\n  },\n  {\n    \"row\": 0,\n    \"col\": 3,\n    \"find\": \"492\u00f75=98, 2\",\n    \"replace\": \"497\u00f72=248, 1\"\n  },\n  {\n    \"row\": 0,\n    \"col\": 4,\n    \"find\": \"739\u00f76=123, 1\",\n    \"replace\": \"305\u00f75=61, 0\"\n  },\n  {\n    \"row\": 4,\n    \"col\": 0,\n    \"find\": \"823\u00f76=137, 1\",\n    \"replace\": \"830\u00f76=138, 2\"\n  },\n  {\n    \"row\": 4,\n    \"col\": 1,\n    \"find\": \"950\u00f73=316, 2\",\n    \"replace\": \"385\u00f76=64, 1\"\n  },\n  {\n    \"row\": 4,\n    \"col\": 2,\n    \"find\": \"426\u00f78=53, 2\",\n    \"replace\": \"349\u00f74=87, 1\"\n  },\n  {\n    \"row\": 4,\n    \"col\": 3,\n    \"find\": \"880\u00f74=220, 0\",\n    \"replace\": \"871\u00f73=290, 1\"\n  },\n  {\n    \"row\": 4,\n    \"col\": 4,\n    \"find\": \"854\u00f73=284, 2\",\n    \"replace\": \"741\u00f78=92, 5\"\n  },\n  {\n    \"row\": 8,\n    \"col\": 0,\n    \"find\": \"102\u00f75=20, 2\",\n    \"replace\": \"640\u00f75=128, 0\"\n  },\n  {\n    \"row\": 8,\n    \"col\": 1,\n    \"find\": \"986\u00f74=246, 2\",\n    \"replace\": \"606\u00f73=202, 0\"\n  },\n  {\n    \"row\": 8,\n    \"col\": 2,\n    \"find\": \"986\u00f76=164, 2\",\n    \"replace\": \"760\u00f73=253, 1\"\n  },\n  {\n    \"row\": 8,\n    \"col\": 3,\n    \"find\": \"964\u00f75=192, 4\",\n    \"replace\": \"920\u00f72=460, 0\"\n  },\n  {\n    \"row\": 8,\n    \"col\": 4,\n    \"find\": \"299\u00f75=59, 4\",\n    \"replace\": \"849\u00f73=283, 0\"\n  },\n  {\n    \"row\": 12,\n    \"col\": 0,\n    \"find\": \"909\u00f72=454, 1\",\n    \"replace\": \"979\u00f72=489, 1\"\n  },\n  {\n    \"row\": 12,\n    \"col\": 1,\n    \"find\": \"931\u00f78=116, 3\",\n    \"replace\": \"648\u00f75=129, 3\"\n  },\n  {\n    \"row\": 12,\n    \"col\": 2,\n    \"find\": \"169\u00f73=56, 1\",\n    \"replace\": \"881\u00f79=97, 8\"\n  },\n  {\n    \"row\": 12,\n    \"col\": 3,\n    \"find\": \"835\u00f78=104, 3\",\n    \"replace\": \"233\u00f79=25, 8\"\n  },\n  {\n    \"row\": 12,\n    \"col\": 4,\n    \"find\": \"273\u00f79=30, 3\",\n    \"replace\": \"482\u00f77=68, 6\"\n  },\n  {\n    \"row\": 16,\n    \"col\": 0,\n    \"find\": \"762\u00f77=108, 6\",\n    \"replace\": \"306\u00f72=153, 0\"\n  },\n  {\n    \"row\": 16,\n    \"col\": 1,\n    \"find\": \"914\u00f76=152, 2\",\n    \"replace\": \"449\u00f79=49, 8\"\n  },\n  {\n    \"row\": 16,\n    \"col\": 2,\n    \"find\": \"776\u00f72=388, 0\",\n    \"replace\": \"190\u00f73=63, 1\"\n  },\n  {\n    \"row\": 16,\n    \"col\": 3,\n    \"find\": \"586\u00f78=73, 2\",\n    \"replace\": \"569\u00f76=94, 5\"\n  },\n  {\n    \"row\": 16,\n    \"col\": 4,\n    \"find\": \"693\u00f72=346, 1\",\n    \"replace\": \"737\u00f79=81, 8\"\n  }\n];\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nif (tables.items.length === 0) {\n  throw new Error(\"Expected a table in the document body, found none.\");\n}\nconst tbl = tables.items[0];\n\n// Load current text for every target cell so we can sanity-check\n// before overwriting (defends against a mismatched/edited source doc).\nconst cells = edits.map((e) => tbl.getCell(e.row, e.col));\ncells.forEach((cell) => cell.body.load(\"text\"));\nawait context.sync();\n\nedits.forEach((e, i) => {\n  const actual = cells[i].body.text.trim();\n  if (actual !== e.find) {\n    throw new Error(\n      `Cell (${e.row}, ${e.col}) text mismatch: expected \"${e.find}\", found \"${actual}\"`\n    );\n  }\n});\n\n// Replace the run's text in place (via the cell body's full range) so\n// the existing run/paragraph formatting (font, size, alignment) is\n// preserved - only the <w:t> contents change.\ncells.forEach((cell, i) => {\n  cell.body.getRange().insertText(edits[i].replace, Word.InsertLocation.replace);\n});\n\nawait context.sync();\n", "ps1": "# Replace each three-digit-division answer cell with its new value.\n# Cells are addressed positionally (row, col) so the edit is immune to\n# collisions between an old value here and a new value placed elsewhere\n# (e.g. \"986\u00f74=246, 2\" is both an old cell value and a new cell value).\n$d = $word.ActiveDocument\n$tbl = $d.Tables.Item(1)\n\n$edits = @(\n  @{ Row=1; Col=1; Find=\"793\u00f75=158, 3\"; Replace=\"986\u00f74=246, 2\" },\n  @{ Row=1; Col=2; Find=\"489\u00f74=122, 1\"; Replace=\"231\u00f77=33, 0\" },\n  @{ Row=1; Col=3; Find=\"371\u00f78=46, 3\"; Replace=\"644\u00f77=92, 0\" },\n  @{ Row=1; Col=4; Find=\"492\u00f75=98, 2\"; Replace=\"497\u00f72=248, 1\" },\n  @{ Row=1; Col=5; Find=\"739\u00f76=123, 1\"; Replace=\"305\u00f75=61, 0\" },\n  @{ Row=5; Col=1; Find=\"823\u00f76=137, 1\"; Replace=\"830\u00f76=138, 2\" },\n  @{ Row=5; Col=2; Find=\"950\u00f73=316, 2\"; Replace=\"385\u00f76=64, 1\" },\n  @{ Row=5; Col=3; Find=\"426\u00f78=53, 2\"; Replace=\"349\u00f74=87, 1\" },\n  @{ Row=5; Col=4; Find=\"880\u00f74=220, 0\"; Replace=\"871\u00f73=290, 1\" },\n  @{ Row=5; Col=5; Find=\"854\u00f73=284, 2\"; Replace=\"741\u00f78=92, 5\" },\n  @{ Row=9; Col=1; Find=\"102\u00f75=20, 2\"; Replace=\"640\u00f75=128, 0\" },\n  @{ Row=9; Col=2; Find=\"986\u00f74=246, 2\"; Replace=\"606\u00f73=202, 0\" },\n  @{ Row=9; Col=3; Find=\"986\u00f76=164, 2\"; Replace=\"760\u00f73=253, 1\" },\n  @{ Row=9; Col=4; Find=\"964\u00f75=192, 4\"; Replace=\"920\u00f72=460, 0\" },\n  @{ Row=9; Col=5; Find=\"299\u00f75=59, 4\"; Replace=\"849\u00f73=283, 0\" },\n  @{ Row=13; Col=1; Find=\"909\u00f72=454, 1\"; Replace=\"979\u00f72=489, 1\" },\n  @{ Row=13; Col=2; Find=\"931\u00f78=116, 3\"; Replace=\"648\u00f75=129, 3\" },\n  @{ Row=13; Col=3; Find=\"169\u00f73=56, 1\"; Replace=\"881\u00f79=97, 8\" },\n  @{ Row=13; Col=4; Find=\"835\u00f78=104, 3\"; Replace=\"233\u00f79=25, 8\" },\n  @{ Row=13; Col=5; Find=\"273\u00f79=30, 3\"; Replace=\"482\u00f77=68, 6\" },\n  @{ Row=17; Col=1; Find=\"762\u00f77=108, 6\"; Replace=\"306\u00f72=153, 0\" },\n  @{ Row=17; Col=2; Find=\"914\u00f76=152, 2\"; Replace=\"449\u00f79=49, 8\" },\n  @{ Row=17; Col=3; Find=\"776\u00f72=388, 0\"; Replace=\"190\u00f73=63, 1\" },\n  @{ Row=17; Col=4; Find=\"586\u00f78=73, 2\"; Replace=\"569\u00f76=94, 5\" },\n  @{ Row=17; Col=5; Find=\"693\u00f72=346, 1\"; Replace=\"737\u00f79=81, 8\" }\n)\n\nforeach ($e in $edits) {\n  $cell = $tbl.Cell($e.Row, $e.Col)\n  $actual = $cell.Range.Text.TrimEnd([char]13, [char]7)\n  if ($actual -ne $e.Find) {\n    throw \"Cell ($($e.Row),$($e.Col)) text mismatch: expected `\"$($e.Find)`\" found `\"$actual`\"\"\n  }\n  $cell.Range.Text = $e.Replace\n}\n"}
